$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Cells.Item(19, 1).Value = 131090310
$ws.Cells.Item(19, 2).Value = 91804
$ws.Cells.Item(19, 4).Value = "NT"
$ws.Cells.Item(19, 5).Value = 1108
$ws.Cells.Item(19, 6).Value = "Harticka"
$ws.Cells.Item(19, 7).Value = "Pelloporus leporinus"
$ws.Cells.Item(19, 8).Value = "(Fr.) Krieglst."
$ws.Cells.Item(19, 16).Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Cells.Item(19, 17).Value = 584977
$ws.Cells.Item(19, 18).Value = 7060194
$ws.Cells.Item(19, 19).Value = 15
$ws.Cells.Item(19, 20).Value = "Västernorrland"
$ws.Cells.Item(19, 21).Value = "Sollefteå"
$ws.Cells.Item(19, 22).Value = "Ångermanland"
$ws.Cells.Item(19, 23).Value = "Junsele"
$ws.Cells.Item(19, 25).Value = "'2026-02-09"
$ws.Cells.Item(19, 26).Value = "14:06"
$ws.Cells.Item(19, 27).Value = "'2026-02-09"
$ws.Cells.Item(19, 28).Value = "14:06"
$ws.Cells.Item(19, 30).Value = $false
$ws.Cells.Item(19, 31).Value = $false
$ws.Cells.Item(19, 33).Value = $false
$ws.Cells.Item(19, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(19, 50).Value = "Daniel Rutschman"

# Row 20
$ws.Cells.Item(20, 1).Value = 131090145
$ws.Cells.Item(20, 2).Value = 79243
$ws.Cells.Item(20, 4).Value = "NT"
$ws.Cells.Item(20, 5).Value = 6425
$ws.Cells.Item(20, 6).Value = "Garnlav"
$ws.Cells.Item(20, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(20, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(20, 16).Value = "Timmeråsen, Timmeråsen, Ång"
$ws.Cells.Item(20, 17).Value = 585013
$ws.Cells.Item(20, 18).Value = 7060142
$ws.Cells.Item(20, 19).Value = 10
$ws.Cells.Item(20, 20).Value = "Västernorrland"
$ws.Cells.Item(20, 21).Value = "Sollefteå"
$ws.Cells.Item(20, 22).Value = "Ångermanland"
$ws.Cells.Item(20, 23).Value = "Junsele"
$ws.Cells.Item(20, 25).Value = "'2026-02-09"
$ws.Cells.Item(20, 26).Value = "13:53"
$ws.Cells.Item(20, 27).Value = "'2026-02-09"
$ws.Cells.Item(20, 28).Value = "13:53"
$ws.Cells.Item(20, 30).Value = $false
$ws.Cells.Item(20, 31).Value = $false
$ws.Cells.Item(20, 33).Value = $false
$ws.Cells.Item(20, 49).Value = "Kim Hultgren"
$ws.Cells.Item(20, 50).Value = "Kim Hultgren"

# Row 21
$ws.Cells.Item(21, 1).Value = 131090091
$ws.Cells.Item(21, 2).Value = 57884
$ws.Cells.Item(21, 4).Value = "NT"
$ws.Cells.Item(21, 5).Value = 100109
$ws.Cells.Item(21, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(21, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(21, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(21, 13).Value = "färska spår"
$ws.Cells.Item(21, 16).Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Cells.Item(21, 17).Value = 585024
$ws.Cells.Item(21, 18).Value = 7060099
$ws.Cells.Item(21, 19).Value = 15
$ws.Cells.Item(21, 20).Value = "Västernorrland"
$ws.Cells.Item(21, 21).Value = "Sollefteå"
$ws.Cells.Item(21, 22).Value = "Ångermanland"
$ws.Cells.Item(21, 23).Value = "Junsele"
$ws.Cells.Item(21, 25).Value = "'2026-02-09"
$ws.Cells.Item(21, 26).Value = "13:50"
$ws.Cells.Item(21, 27).Value = "'2026-02-09"
$ws.Cells.Item(21, 28).Value = "13:50"
$ws.Cells.Item(21, 29).Value = "Färska ringhack, tall"
$ws.Cells.Item(21, 30).Value = $false
$ws.Cells.Item(21, 31).Value = $false
$ws.Cells.Item(21, 33).Value = $false
$ws.Cells.Item(21, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(21, 50).Value = "Daniel Rutschman"

# Row 22
$ws.Cells.Item(22, 1).Value = 131090374
$ws.Cells.Item(22, 2).Value = 57884
$ws.Cells.Item(22, 4).Value = "NT"
$ws.Cells.Item(22, 5).Value = 100109
$ws.Cells.Item(22, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(22, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(22, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(22, 13).Value = "färska spår"
$ws.Cells.Item(22, 16).Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Cells.Item(22, 17).Value = 584977
$ws.Cells.Item(22, 18).Value = 7060203
$ws.Cells.Item(22, 19).Value = 15
$ws.Cells.Item(22, 20).Value = "Västernorrland"
$ws.Cells.Item(22, 21).Value = "Sollefteå"
$ws.Cells.Item(22, 22).Value = "Ångermanland"
$ws.Cells.Item(22, 23).Value = "Junsele"
$ws.Cells.Item(22, 25).Value = "'2026-02-09"
$ws.Cells.Item(22, 26).Value = "14:08"
$ws.Cells.Item(22, 27).Value = "'2026-02-09"
$ws.Cells.Item(22, 28).Value = "14:08"
$ws.Cells.Item(22, 29).Value = "Färska ringhack tall"
$ws.Cells.Item(22, 30).Value = $false
$ws.Cells.Item(22, 31).Value = $false
$ws.Cells.Item(22, 33).Value = $false
$ws.Cells.Item(22, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(22, 50).Value = "Daniel Rutschman"

# Row 23
$ws.Cells.Item(23, 1).Value = 131090601
$ws.Cells.Item(23, 2).Value = 57884
$ws.Cells.Item(23, 4).Value = "NT"
$ws.Cells.Item(23, 5).Value = 100109
$ws.Cells.Item(23, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(23, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(23, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(23, 13).Value = "färska spår"
$ws.Cells.Item(23, 16).Value = "Hållflon, Hållflon, Ång"
$ws.Cells.Item(23, 17).Value = 584871
$ws.Cells.Item(23, 18).Value = 7060419
$ws.Cells.Item(23, 19).Value = 10
$ws.Cells.Item(23, 20).Value = "Västernorrland"
$ws.Cells.Item(23, 21).Value = "Sollefteå"
$ws.Cells.Item(23, 22).Value = "Ångermanland"
$ws.Cells.Item(23, 23).Value = "Junsele"
$ws.Cells.Item(23, 25).Value = "'2026-02-09"
$ws.Cells.Item(23, 26).Value = "14:30"
$ws.Cells.Item(23, 27).Value = "'2026-02-09"
$ws.Cells.Item(23, 28).Value = "14:30"
$ws.Cells.Item(23, 29).Value = "Ringhack på tall"
$ws.Cells.Item(23, 30).Value = $false
$ws.Cells.Item(23, 31).Value = $false
$ws.Cells.Item(23, 33).Value = $false
$ws.Cells.Item(23, 49).Value = "Kim Hultgren"
$ws.Cells.Item(23, 50).Value = "Kim Hultgren"

# Row 24
$ws.Cells.Item(24, 1).Value = 131090313
$ws.Cells.Item(24, 2).Value = 57884
$ws.Cells.Item(24, 4).Value = "NT"
$ws.Cells.Item(24, 5).Value = 100109
$ws.Cells.Item(24, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(24, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(24, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(24, 13).Value = "färska spår"
$ws.Cells.Item(24, 16).Value = "Timmeråsen, Timmeråsen, Ång"
$ws.Cells.Item(24, 17).Value = 585012
$ws.Cells.Item(24, 18).Value = 7060182
$ws.Cells.Item(24, 19).Value = 10
$ws.Cells.Item(24, 20).Value = "Västernorrland"
$ws.Cells.Item(24, 21).Value = "Sollefteå"
$ws.Cells.Item(24, 22).Value = "Ångermanland"
$ws.Cells.Item(24, 23).Value = "Junsele"
$ws.Cells.Item(24, 25).Value = "'2026-02-09"
$ws.Cells.Item(24, 26).Value = "14:05"
$ws.Cells.Item(24, 27).Value = "'2026-02-09"
$ws.Cells.Item(24, 28).Value = "14:05"
$ws.Cells.Item(24, 29).Value = "Ringhack på tall"
$ws.Cells.Item(24, 30).Value = $false
$ws.Cells.Item(24, 31).Value = $false
$ws.Cells.Item(24, 33).Value = $false
$ws.Cells.Item(24, 49).Value = "Kim Hultgren"
$ws.Cells.Item(24, 50).Value = "Kim Hultgren"

# Row 25
$ws.Cells.Item(25, 1).Value = 131090020
$ws.Cells.Item(25, 2).Value = 57884
$ws.Cells.Item(25, 4).Value = "NT"
$ws.Cells.Item(25, 5).Value = 100109
$ws.Cells.Item(25, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(25, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(25, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(25, 13).Value = "färska spår"
$ws.Cells.Item(25, 16).Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Cells.Item(25, 17).Value = 585065
$ws.Cells.Item(25, 18).Value = 7060076
$ws.Cells.Item(25, 19).Value = 15
$ws.Cells.Item(25, 20).Value = "Västernorrland"
$ws.Cells.Item(25, 21).Value = "Sollefteå"
$ws.Cells.Item(25, 22).Value = "Ångermanland"
$ws.Cells.Item(25, 23).Value = "Junsele"
$ws.Cells.Item(25, 25).Value = "'2026-02-09"
$ws.Cells.Item(25, 26).Value = "13:47"
$ws.Cells.Item(25, 27).Value = "'2026-02-09"
$ws.Cells.Item(25, 28).Value = "13:47"
$ws.Cells.Item(25, 29).Value = "Färska ringhack tall"
$ws.Cells.Item(25, 30).Value = $false
$ws.Cells.Item(25, 31).Value = $false
$ws.Cells.Item(25, 33).Value = $false
$ws.Cells.Item(25, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(25, 50).Value = "Daniel Rutschman"

# Row 26
$ws.Cells.Item(26, 1).Value = 131090008
$ws.Cells.Item(26, 2).Value = 57884
$ws.Cells.Item(26, 4).Value = "NT"
$ws.Cells.Item(26, 5).Value = 100109
$ws.Cells.Item(26, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(26, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(26, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(26, 13).Value = "färska spår"
$ws.Cells.Item(26, 16).Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Cells.Item(26, 17).Value = 585080
$ws.Cells.Item(26, 18).Value = 7060076
$ws.Cells.Item(26, 19).Value = 15
$ws.Cells.Item(26, 20).Value = "Västernorrland"
$ws.Cells.Item(26, 21).Value = "Sollefteå"
$ws.Cells.Item(26, 22).Value = "Ångermanland"
$ws.Cells.Item(26, 23).Value = "Junsele"
$ws.Cells.Item(26, 25).Value = "'2026-02-09"
$ws.Cells.Item(26, 26).Value = "13:46"
$ws.Cells.Item(26, 27).Value = "'2026-02-09"
$ws.Cells.Item(26, 28).Value = "13:46"
$ws.Cells.Item(26, 29).Value = "Färska ringhack, tall"
$ws.Cells.Item(26, 30).Value = $false
$ws.Cells.Item(26, 31).Value = $false
$ws.Cells.Item(26, 33).Value = $false
$ws.Cells.Item(26, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(26, 50).Value = "Daniel Rutschman"

# Row 27
$ws.Cells.Item(27, 1).Value = 131090015
$ws.Cells.Item(27, 2).Value = 57884
$ws.Cells.Item(27, 4).Value = "NT"
$ws.Cells.Item(27, 5).Value = 100109
$ws.Cells.Item(27, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(27, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(27, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(27, 13).Value = "färska spår"
$ws.Cells.Item(27, 16).Value = "Timmeråsen, Timmeråsen, Ång"
$ws.Cells.Item(27, 17).Value = 585076
$ws.Cells.Item(27, 18).Value = 7060075
$ws.Cells.Item(27, 19).Value = 10
$ws.Cells.Item(27, 20).Value = "Västernorrland"
$ws.Cells.Item(27, 21).Value = "Sollefteå"
$ws.Cells.Item(27, 22).Value = "Ångermanland"
$ws.Cells.Item(27, 23).Value = "Junsele"
$ws.Cells.Item(27, 25).Value = "'2026-02-09"
$ws.Cells.Item(27, 26).Value = "13:46"
$ws.Cells.Item(27, 27).Value = "'2026-02-09"
$ws.Cells.Item(27, 28).Value = "13:46"
$ws.Cells.Item(27, 29).Value = "Ringhack på tall"
$ws.Cells.Item(27, 30).Value = $false
$ws.Cells.Item(27, 31).Value = $false
$ws.Cells.Item(27, 33).Value = $false
$ws.Cells.Item(27, 49).Value = "Kim Hultgren"
$ws.Cells.Item(27, 50).Value = "Kim Hultgren"

# Row 28
$ws.Cells.Item(28, 1).Value = 131089521
$ws.Cells.Item(28, 2).Value = 57881
$ws.Cells.Item(28, 4).Value = "NT"
$ws.Cells.Item(28, 5).Value = 100049
$ws.Cells.Item(28, 6).Value = "Spillkråka"
$ws.Cells.Item(28, 7).Value = "Dryocopus martius"
$ws.Cells.Item(28, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(28, 13).Value = "färska spår"
$ws.Cells.Item(28, 16).Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Cells.Item(28, 17).Value = 584995
$ws.Cells.Item(28, 18).Value = 7060537
$ws.Cells.Item(28, 19).Value = 15
$ws.Cells.Item(28, 20).Value = "Västernorrland"
$ws.Cells.Item(28, 21).Value = "Sollefteå"
$ws.Cells.Item(28, 22).Value = "Ångermanland"
$ws.Cells.Item(28, 23).Value = "Junsele"
$ws.Cells.Item(28, 25).Value = "'2026-02-09"
$ws.Cells.Item(28, 26).Value = "13:14"
$ws.Cells.Item(28, 27).Value = "'2026-02-09"
$ws.Cells.Item(28, 28).Value = "13:14"
$ws.Cells.Item(28, 30).Value = $false
$ws.Cells.Item(28, 31).Value = $false
$ws.Cells.Item(28, 33).Value = $false
$ws.Cells.Item(28, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(28, 50).Value = "Daniel Rutschman"

# Row 29
$ws.Cells.Item(29, 1).Value = 131090275
$ws.Cells.Item(29, 2).Value = 57884
$ws.Cells.Item(29, 4).Value = "NT"
$ws.Cells.Item(29, 5).Value = 100109
$ws.Cells.Item(29, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(29, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(29, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(29, 12).Value = "hane"
$ws.Cells.Item(29, 13).Value = "födosökande"
$ws.Cells.Item(29, 16).Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Cells.Item(29, 17).Value = 584987
$ws.Cells.Item(29, 18).Value = 7060190
$ws.Cells.Item(29, 19).Value = 15
$ws.Cells.Item(29, 20).Value = "Västernorrland"
$ws.Cells.Item(29, 21).Value = "Sollefteå"
$ws.Cells.Item(29, 22).Value = "Ångermanland"
$ws.Cells.Item(29, 23).Value = "Junsele"
$ws.Cells.Item(29, 25).Value = "'2026-02-09"
$ws.Cells.Item(29, 27).Value = "'2026-02-09"
$ws.Cells.Item(29, 30).Value = $false
$ws.Cells.Item(29, 31).Value = $false
$ws.Cells.Item(29, 33).Value = $false
$ws.Cells.Item(29, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(29, 50).Value = "Daniel Rutschman"
